$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row 1 (existing data stays in rows 2-17, not shifted)
$ws.Range("A1").Value = "BSA µg/µL"
$ws.Range("B1").Value = "BSA Absorbances"
$ws.Range("C1").Value = "BSA Absorbances Duplicate"
$ws.Range("D1").Value = "Sample Names"
$ws.Range("E1").Value = "Sample Absorbances"
$ws.Range("F1").Value = "Sample Absorbances Duplicate"

# Set column widths per the diff
$ws.Columns.Item(2).ColumnWidth = 16.36328125
$ws.Columns.Item(3).ColumnWidth = 23.453125
$ws.Columns.Item(4).ColumnWidth = 13.08984375
$ws.Columns.Item(5).ColumnWidth = 18.1796875

# Update selection to H8
$ws.Range("H8").Select()
